$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4
$ws.Range("A2").Value = 6
$ws.Range("A3").Value = 7
$ws.Range("A4").Value = 9

$ws.Range("A5").Select()
